# daily updates Sept. 11
# Adds the Sept-10 (row 72, date 45544) and Sept-11 (row 73, date 45545)
# daily counts to the "Babine" sheet and extends the running-total
# formulas (M, N, O, P, Q, R) down through those two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Babine")

# ---- Row 72 (2024-09-10) -------------------------------------------------
$ws.Range("B72").Value = 5248
$ws.Range("C72").Value = 450
$ws.Range("D72").Value = 260
$ws.Range("E72").Value = 6580
$ws.Range("F72").Value = 78
$ws.Range("G72").Value = 39
$ws.Range("H72").Value = 1
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0

$ws.Range("M72").Formula = "=M71+B72"
$ws.Range("N72").Formula = "=J72+L72"
$ws.Range("O72").Formula = "=B72+I72+N72"
$ws.Range("P72").Formula = "=P71+O72"
$ws.Range("Q72").Formula = "=C72+K72"
$ws.Range("R72").Formula = "=Q72+R71"

# ---- Row 73 (2024-09-11) -------------------------------------------------
$ws.Range("B73").Value = 3708
$ws.Range("C73").Value = 448
$ws.Range("D73").Value = 367
$ws.Range("E73").Value = 6326
$ws.Range("F73").Value = 26
$ws.Range("G73").Value = 29
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0

$ws.Range("M73").Formula = "=M72+B73"
$ws.Range("N73").Formula = "=J73+L73"
$ws.Range("O73").Formula = "=B73+I73+N73"
$ws.Range("P73").Formula = "=P72+O73"
$ws.Range("Q73").Formula = "=C73+K73"
$ws.Range("R73").Formula = "=Q73+R72"

# ---- View state: move the active selection to Q70 ------------------------
[void]$ws.Range("Q70").Select()
